$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format risky numeric-looking price cells as Text (one
# statement per cell) so COM's automatic type inference doesn't
# convert them to numbers; the source stores these as literal text.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "71.539.93"
$ws.Range("E2").Value = "  +2.92%  "
$ws.Range("D3").Value = "4.027.20"
$ws.Range("E3").Value = "  +3.08%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "517.08"
$ws.Range("E5").Value = "  -2.29%  "
$ws.Range("D6").Value = "146.75"
$ws.Range("E6").Value = "  +1.49%  "
$ws.Range("D7").Value = "0.618"
$ws.Range("E7").Value = "  +0.80%  "
$ws.Range("E8").Value = "  +0.19%  "
$ws.Range("D9").Value = "0.731"
$ws.Range("E9").Value = "  +0.44%  "
$ws.Range("D10").Value = "0.172"
$ws.Range("E10").Value = "  +0.69%  "
$ws.Range("D11").Value = "0.0000332"
$ws.Range("E11").Value = "  +0.71%  "
$ws.Range("D12").Value = "46.40"
$ws.Range("E12").Value = "  +9.72%  "
$ws.Range("D13").Value = "10.74"
$ws.Range("E13").Value = "  +4.90%  "
$ws.Range("D14").Value = "4.690.65"
$ws.Range("E14").Value = "  +3.30%  "
$ws.Range("D15").Value = "4.046.29"
$ws.Range("E15").Value = "  +3.41%  "
$ws.Range("D16").Value = "21.11"
$ws.Range("E16").Value = "  +6.44%  "
$ws.Range("D17").Value = "14.15"
$ws.Range("E17").Value = "  +2.17%  "
$ws.Range("E18").Value = "  -2.55%  "
$ws.Range("E19").Value = "  -2.10%  "
$ws.Range("D20").Value = "71.692.15"
$ws.Range("E20").Value = "  +3.29%  "
$ws.Range("D21").Value = "434.86"
$ws.Range("E21").Value = "  +0.29%  "
$ws.Range("D22").Value = "95.02"
$ws.Range("E22").Value = "  +8.36%  "
$ws.Range("D23").Value = "3.49"
$ws.Range("E23").Value = "  +4.56%  "
$ws.Range("D24").Value = "14.33"
$ws.Range("E24").Value = "  +0.71%  "
$ws.Range("D25").Value = "11.98"
$ws.Range("E25").Value = "  +3.62%  "
$ws.Range("D26").Value = "4.03"
$ws.Range("E26").Value = "  -0.99%  "
$ws.Range("D27").Value = "11.16"
$ws.Range("E27").Value = "  +4.98%  "
$ws.Range("D28").Value = "36.71"
$ws.Range("E28").Value = "  +1.74%  "
$ws.Range("D29").Value = "3.06"
$ws.Range("E29").Value = "  +8.29%  "
$ws.Range("D30").Value = "703.41"
$ws.Range("E30").Value = "  +1.80%  "
$ws.Range("D31").Value = "13.37"
$ws.Range("E31").Value = "  +1.73%  "
$ws.Range("D32").Value = "0.129"
$ws.Range("E32").Value = "  +2.43%  "
$ws.Range("D33").Value = "7.07"
$ws.Range("E33").Value = "  +20.34%  "
$ws.Range("D34").Value = "67.60"
$ws.Range("E34").Value = "  -0.30%  "
$ws.Range("D35").Value = "0.0₃0893"
$ws.Range("E35").Value = "  +7.49%  "

# Row 36 - was TheGraph, now ThetaToken (rows 36/37 swap rank order)
$ws.Range("B36").Value = "ThetaToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D36").Value = "3.74"
$ws.Range("E36").Value = "  +27.67%  "

# Row 37 - was ThetaToken, now TheGraph
$ws.Range("B37").Value = "TheGraph"
$ws.Range("C37").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D37").Value = "0.439"
$ws.Range("E37").Value = "  -0.80%  "

$ws.Range("D38").Value = "40.37"
$ws.Range("E38").Value = "  +0.56%  "
$ws.Range("D39").Value = "0.153"
$ws.Range("E39").Value = "  +2.44%  "
$ws.Range("E40").Value = "  +0.06%  "
$ws.Range("E41").Value = "  -0.08%  "
$ws.Range("D42").Value = "0.0482"
$ws.Range("E42").Value = "  +0.59%  "
$ws.Range("E43").Value = "  +1.37%  "
$ws.Range("D44").Value = "2.76"
$ws.Range("E44").Value = "  +0.58%  "
$ws.Range("E45").Value = "  +4.18%  "
$ws.Range("D46").Value = "0.144"
$ws.Range("E46").Value = "  +2.20%  "
$ws.Range("D47").Value = "3.14"
$ws.Range("E47").Value = "  +2.88%  "
$ws.Range("D48").Value = "0.000275"
$ws.Range("E48").Value = "  +22.26%  "
$ws.Range("D49").Value = "8.97"
$ws.Range("E49").Value = "  +6.05%  "
$ws.Range("D50").Value = "3.30"
$ws.Range("E50").Value = "  +0.22%  "
$ws.Range("D51").Value = "0.0₆0339"
$ws.Range("E51").Value = "  +0.29%  "
